$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the expected title text in D4 from "Invalid username and password"
# to "Invalid username or password"
$ws.Range("D4").Value = "Invalid username or password"

# Move the active selection to D4 (last edited cell)
$ws.Range("D4").Select()
